$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 9.295069333333332
$ws.Range("H2").Value = 27.885208
$ws.Range("I2").Value = 0.2851098797714356
$ws.Range("J2").Value = 0.2851098797714357
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.233575666666667
$ws.Range("N2").Value = 12.700727
$ws.Range("O2").Value = 0.4702904532725788
$ws.Range("P2").Value = 0.4702904532725788
$ws.Range("Q2").Value = 39.35137934957955
$ws.Range("R2").Value = 354.162414146216
$ws.Range("S2").Value = 0.1340844545901989
$ws.Range("T2").Value = 0.1340844545901989

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 9.295069333333332
$ws.Range("H3").Value = 27.885208
$ws.Range("I3").Value = 0.2851098797714356
$ws.Range("J3").Value = 0.2851098797714357
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.605649999999999
$ws.Range("N3").Value = 10.81695
$ws.Range("O3").Value = 0.4005367817548413
$ws.Range("P3").Value = 0.4005367817548413
$ws.Range("Q3").Value = 33.51476674173333
$ws.Range("R3").Value = 301.6329006756
$ws.Range("S3").Value = 0.1141969936901606
$ws.Range("T3").Value = 0.1141969936901606

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 9.295069333333332
$ws.Range("H4").Value = 27.885208
$ws.Range("I4").Value = 0.2851098797714356
$ws.Range("J4").Value = 0.2851098797714357
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.162819
$ws.Range("N4").Value = 3.488456999999999
$ws.Range("O4").Value = 0.1291727649725799
$ws.Range("P4").Value = 0.1291727649725799
$ws.Range("Q4").Value = 10.80848322711733
$ws.Range("R4").Value = 97.27634904405598
$ws.Range("S4").Value = 0.03682843149107617
$ws.Range("T4").Value = 0.03682843149107618

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 8.588082333333332
$ws.Range("H5").Value = 25.764247
$ws.Range("I5").Value = 0.2634242988100204
$ws.Range("J5").Value = 0.2634242988100204
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.233575666666667
$ws.Range("N5").Value = 12.700727
$ws.Range("O5").Value = 0.4702904532725788
$ws.Range("P5").Value = 0.4702904532725788
$ws.Range("Q5").Value = 36.35829638972989
$ws.Range("R5").Value = 327.224667507569
$ws.Range("S5").Value = 0.1238859328903757
$ws.Range("T5").Value = 0.1238859328903757

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 8.588082333333332
$ws.Range("H6").Value = 25.764247
$ws.Range("I6").Value = 0.2634242988100204
$ws.Range("J6").Value = 0.2634242988100204
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.605649999999999
$ws.Range("N6").Value = 10.81695
$ws.Range("O6").Value = 0.4005367817548413
$ws.Range("P6").Value = 0.4005367817548413
$ws.Range("Q6").Value = 30.96561906518333
$ws.Range("R6").Value = 278.6905715866499
$ws.Range("S6").Value = 0.1055111208813913
$ws.Range("T6").Value = 0.1055111208813913

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 8.588082333333332
$ws.Range("H7").Value = 25.764247
$ws.Range("I7").Value = 0.2634242988100204
$ws.Range("J7").Value = 0.2634242988100204
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.162819
$ws.Range("N7").Value = 3.488456999999999
$ws.Range("O7").Value = 0.1291727649725799
$ws.Range("P7").Value = 0.1291727649725799
$ws.Range("Q7").Value = 9.98638531076433
$ws.Range("R7").Value = 89.87746779687897
$ws.Range("S7").Value = 0.03402724503825344
$ws.Range("T7").Value = 0.03402724503825344

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.71855733333333
$ws.Range("H8").Value = 44.155672
$ws.Range("I8").Value = 0.4514658214185439
$ws.Range("J8").Value = 0.4514658214185439
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.233575666666667
$ws.Range("N8").Value = 12.700727
$ws.Range("O8").Value = 0.4702904532725788
$ws.Range("P8").Value = 0.4702904532725788
$ws.Range("Q8").Value = 62.31212617483823
$ws.Range("R8").Value = 560.809135573544
$ws.Range("S8").Value = 0.2123200657920041
$ws.Range("T8").Value = 0.2123200657920041

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.71855733333333
$ws.Range("H9").Value = 44.155672
$ws.Range("I9").Value = 0.4514658214185439
$ws.Range("J9").Value = 0.4514658214185439
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.605649999999999
$ws.Range("N9").Value = 10.81695
$ws.Range("O9").Value = 0.4005367817548413
$ws.Range("P9").Value = 0.4005367817548413
$ws.Range("Q9").Value = 53.06996624893333
$ws.Range("R9").Value = 477.6296962403999
$ws.Range("S9").Value = 0.1808286671832895
$ws.Range("T9").Value = 0.1808286671832894

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 14.71855733333333
$ws.Range("H10").Value = 44.155672
$ws.Range("I10").Value = 0.4514658214185439
$ws.Range("J10").Value = 0.4514658214185439
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.162819
$ws.Range("N10").Value = 3.488456999999999
$ws.Range("O10").Value = 0.1291727649725799
$ws.Range("P10").Value = 0.1291727649725799
$ws.Range("Q10").Value = 17.11501811978933
$ws.Range("R10").Value = 154.035163078104
$ws.Range("S10").Value = 0.05831708844325031
$ws.Range("T10").Value = 0.05831708844325031
